# Handback status report generation: refresh the "latest generate" /
# "correspond handback datetime" timestamps for the first (ccaf072b...)
# file row on each sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-31 04:50:01"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-31 04:49:54"
$wsZhCn.Range("K2").Value = "2016-08-31 04:50:22"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-31 04:50:01"
$wsDeDe.Range("K2").Value = "2016-08-31 04:50:29"
